$wb = $excel.ActiveWorkbook

# New row of data (time = 2025-06-16) to be appended as row 38 on every
# worksheet. Each sheet repeats the same values as its row 37, except for
# the "time" column (A), which advances by one day for all sheets.
$newTime = 45824.43559027778

$rowsData = @{
    1 = @("0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x01,0x70", "0x14", 380, "7.598631275147109e+23", 368, 14)
    2 = @("0x01,0x7c", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x01,0x70", "0xe",  380, "5.68432987514711e+23",  368, 14)
    3 = @("0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x80", "0x7",  130, "5.68631262647114e+23",  128, 7)
    4 = @("0x00,0x82", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x80", "0x3",  130, "9.85046333984776e+23",  128, 3)
}

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $vals = $rowsData[$i]

    $ws.Cells.Item(38, 1).Value = $newTime
    $ws.Cells.Item(38, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item(38, 2).Value = $vals[0]
    $ws.Cells.Item(38, 3).Value = $vals[1]
    $ws.Cells.Item(38, 4).Value = $vals[2]
    $ws.Cells.Item(38, 5).Value = $vals[3]
    $ws.Cells.Item(38, 6).Value = $vals[4]
    $ws.Cells.Item(38, 7).Value = [double]$vals[5]
    $ws.Cells.Item(38, 8).Value = $vals[6]
    $ws.Cells.Item(38, 9).Value = $vals[7]
}
